# Regenerate merged AHB files
# 1. Rename header labels: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
# 2. Turn the data range into an Excel Table (ListObject)
# 3. Freeze the header row (split below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row labels -------------------------------------
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old -replace '_old$', '_FV2310'
        $new = $new -replace '_new$', '_FV2404'
        if ($new -ne $old) {
            $cell.Value = $new
        }
    }
}

# --- 2) Convert the data range A1:U71 into an Excel Table -----------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3) Freeze panes below the header row ---------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
